# Auto-generated edit script
# Applies 2026-02-09 violent crime data update across 37 worksheets (139 cell changes)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("M2").Value = 541
$ws.Range("L3").Value = 7109
$ws.Range("M3").Value = 601
$ws.Range("H4").Value = 1772
$ws.Range("I4").Value = 1849
$ws.Range("M4").Value = 177
$ws.Range("M6").Value = 467
$ws.Range("H7").Value = 26087
$ws.Range("I7").Value = 26319
$ws.Range("L7").Value = 21733
$ws.Range("M7").Value = 1824

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("M2").Value = 37
$ws.Range("M6").Value = 35
$ws.Range("M7").Value = 122

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("M3").Value = 21
$ws.Range("M7").Value = 62

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("M6").Value = 22
$ws.Range("M7").Value = 76

$ws = $wb.Worksheets.Item('New City')
$ws.Range("M2").Value = 17
$ws.Range("M7").Value = 44

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("M6").Value = 4
$ws.Range("M7").Value = 6

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("M7").Value = 54
$ws.Range("M8").Value = 122
$ws.Range("M11").Value = 24
$ws.Range("M18").Value = 17
$ws.Range("M19").Value = 61
$ws.Range("M20").Value = 56
$ws.Range("M21").Value = 5
$ws.Range("M24").Value = 9
$ws.Range("M29").Value = 93
$ws.Range("M30").Value = 6
$ws.Range("M31").Value = 21
$ws.Range("M33").Value = 62
$ws.Range("M36").Value = 19
$ws.Range("M37").Value = 76
$ws.Range("M42").Value = 62
$ws.Range("M43").Value = 16
$ws.Range("M47").Value = 15
$ws.Range("M48").Value = 26
$ws.Range("M50").Value = 9
$ws.Range("L51").Value = 268
$ws.Range("M51").Value = 26
$ws.Range("M52").Value = 24
$ws.Range("M54").Value = 32
$ws.Range("M55").Value = 19
$ws.Range("H63").Value = 323
$ws.Range("I63").Value = 270
$ws.Range("M63").Value = 4
$ws.Range("M65").Value = 44
$ws.Range("M67").Value = 53
$ws.Range("M76").Value = 27
$ws.Range("M77").Value = 17
$ws.Range("M84").Value = 10
$ws.Range("L85").Value = 1085
$ws.Range("M85").Value = 98
$ws.Range("M86").Value = 13
$ws.Range("L88").Value = 226
$ws.Range("M89").Value = 26
$ws.Range("M91").Value = 27
$ws.Range("M94").Value = 23
$ws.Range("M97").Value = 17
$ws.Range("H101").Value = 26087
$ws.Range("I101").Value = 26319
$ws.Range("L101").Value = 21733
$ws.Range("M101").Value = 1824

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("M2").Value = 7
$ws.Range("M6").Value = 6
$ws.Range("M7").Value = 21

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("M3").Value = 15
$ws.Range("M7").Value = 53

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("M6").Value = 4
$ws.Range("M7").Value = 10

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("M4").Value = 4
$ws.Range("M7").Value = 32

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("M2").Value = 28
$ws.Range("M3").Value = 34
$ws.Range("M6").Value = 23
$ws.Range("M7").Value = 93

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("M6").Value = 15
$ws.Range("M7").Value = 26

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("M3").Value = 22
$ws.Range("M6").Value = 17
$ws.Range("M7").Value = 61

$ws = $wb.Worksheets.Item('River North')
$ws.Range("M2").Value = 6
$ws.Range("M3").Value = 6
$ws.Range("M6").Value = 12
$ws.Range("M7").Value = 27

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("M3").Value = 22
$ws.Range("M7").Value = 62

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("M3").Value = 6
$ws.Range("M7").Value = 19

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("M3").Value = 5
$ws.Range("M7").Value = 9

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("M2").Value = 6
$ws.Range("M4").Value = 3
$ws.Range("M7").Value = 27

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("M6").Value = 1
$ws.Range("M7").Value = 5

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("M2").Value = 22
$ws.Range("M7").Value = 56

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("M4").Value = 1
$ws.Range("M7").Value = 17

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("M6").Value = 4
$ws.Range("M7").Value = 19

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("M3").Value = 18
$ws.Range("M7").Value = 54

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("M3").Value = 6
$ws.Range("M7").Value = 23

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("M3").Value = 11
$ws.Range("M7").Value = 15

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("M2").Value = 2
$ws.Range("M7").Value = 9

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("M6").Value = 6
$ws.Range("M7").Value = 24

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("M2").Value = 5
$ws.Range("M7").Value = 17

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L3").Value = 76
$ws.Range("L7").Value = 226

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("M6").Value = 8
$ws.Range("M7").Value = 26

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("M4").Value = 8
$ws.Range("M7").Value = 13

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L4").Value = 41
$ws.Range("M4").Value = 5
$ws.Range("L7").Value = 268
$ws.Range("M7").Value = 26

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("M4").Value = 3
$ws.Range("M7").Value = 16

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("M2").Value = 27
$ws.Range("M3").Value = 43
$ws.Range("L4").Value = 65
$ws.Range("M6").Value = 20
$ws.Range("L7").Value = 1085
$ws.Range("M7").Value = 98

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("M3").Value = 7
$ws.Range("M7").Value = 17

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("M4").Value = 2
$ws.Range("M7").Value = 24

Write-Host "Applied 139 cell updates across 37 worksheets"